# Corrects the path of the projectile model shown in the "Try Out the
# Code" slide's terminal textbox, and fixes the "visualstudio code" typo
# on the "Your First RAVEN Input" slide.

$p = $ppt.ActivePresentation

# --- Slide 5: "Try Out the Code" -> TextBox 3 (the terminal transcript) ---
$slide5 = $p.Slides.Item(5)
$textBox = $slide5.Shapes.Item(3)
$tr = $textBox.TextFrame.TextRange

# Work right-to-left across the "> cd doc/workshop/forwardSampling/projectile_model"
# line so earlier character offsets (computed against the original text)
# stay valid as later edits change the string length.

# Remove the "projectile_model" run entirely.
$tr.Characters(69, 16).Text = ""

# Remove the "/" run entirely.
$tr.Characters(68, 1).Text = ""

# "forwardSampling" -> "ExternalModels" (keeps that run's own formatting).
$tr.Characters(53, 15).Text = "ExternalModels"

# " cd doc/workshop/" -> " cd ~/projects/raven/doc/workshop/" (keeps that
# run's own formatting).
$tr.Characters(36, 17).Text = " cd ~/projects/raven/doc/workshop/"

# Remove the trailing blank paragraph after the "-o out" line; the text
# box auto-fits, so its shape height re-derives from the new text extent.
$tr.Paragraphs(5, 1).Delete()

# --- Slide 9: "1) Your First RAVEN Input" -> editor bullet list ---
$slide9 = $p.Slides.Item(9)
$contentPlaceholder = $slide9.Shapes.Item(2)
$tr9 = $contentPlaceholder.TextFrame.TextRange

# Expand " code" into "Visual Studio Code" first (this keeps the second
# run's formatting, which lacks the err="1" spell-check flag), then drop
# the now-redundant leading "visualstudio" run.
$tr9.Characters(96, 5).Text = "Visual Studio Code"
$tr9.Characters(84, 12).Text = ""
